$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so decimal-like strings (e.g. "1.0000")
# are stored verbatim instead of being auto-coerced to numbers by Excel.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

# --- Cell value updates (per commit diff) ---
$ws.Range("D2").Value = '29.269.52'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '1.871.52'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '0.7118'
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").Value = '241.43'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.3114'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.07667'
$ws.Range("E9").Value = '  -3.46%  '
$ws.Range("D10").Value = '24.72'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("D11").Value = '0.08407'
$ws.Range("E11").Value = '  +1.66%  '
$ws.Range("D12").Value = '1.890.95'
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").Value = '5.234'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").Value = '0.7113'
$ws.Range("E14").Value = '  -2.92%  '
$ws.Range("D15").Value = '91.34'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '29.279.42'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").Value = '5.940'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").Value = '243.62'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").Value = '0.000007865'
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").Value = '2.117.01'
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").Value = '13.15'
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '7.864'
$ws.Range("E23").Value = '  -2.43%  '
$ws.Range("D24").Value = '0.9996'
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").Value = '0.1644'
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("D26").Value = '163.77'
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").Value = '8.995'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '18.52'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = '1.508'
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("D30").Value = '1.309'
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("D31").Value = '4.398'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").Value = '4.259'
$ws.Range("E32").Value = '  +3.44%  '
$ws.Range("D33").Value = '0.05163'
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '1.915'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7730'
$ws.Range("E35").Value = '  +6.32%  '
$ws.Range("D36").Value = '1.170'
$ws.Range("E36").Value = '  -2.52%  '
$ws.Range("D37").Value = '2.684'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '0.01858'
$ws.Range("E38").Value = '  -0.97%  '
$ws.Range("D39").Value = '2.710'
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.160.57'
$ws.Range("E40").Value = '  -3.75%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.400'
$ws.Range("E41").Value = '  +4.09%  '
$ws.Range("D42").Value = '0.8924'
$ws.Range("E42").Value = '  -2.15%  '
$ws.Range("D43").Value = '73.30'
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '103.86'
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '0.9994'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Value = '2.014.57'
$ws.Range("D47").Value = '0.5175'
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("D48").Value = '1.785'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("D49").Value = '9.386'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("D51").Value = '0.4299'
$ws.Range("E51").Value = '  -0.76%  '

# Restore column D to the default "Normal" style (matches original,
# unstyled cells) now that the text values are safely stored.
$priceCol.Style = "Normal"

